$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: fill in the previously-empty "Khuyen mai" (promo) cell
$ws.Cells.Item(2, 5).Value = "KM001"

# New rows 6-9 contain the additional invoices.
# Column B holds date-like text ("2023-12-09"); pre-format those cells as
# Text so Excel keeps them as strings instead of converting to date serials,
# matching the rest of the sheet (which stores dates as plain text).
$ws.Range("B6:B9").NumberFormat = "@"

# Row 6
$ws.Cells.Item(6, 1).Value = "HD005"
$ws.Cells.Item(6, 2).Value = "2023-12-09"
$ws.Cells.Item(6, 3).Value = "KH041"
$ws.Cells.Item(6, 4).Value = "NV003"
$ws.Cells.Item(6, 5).Value = "KM001,"
$ws.Cells.Item(6, 6).Value = 298000.0

# Row 7
$ws.Cells.Item(7, 1).Value = "HD006"
$ws.Cells.Item(7, 2).Value = "2023-12-09"
$ws.Cells.Item(7, 3).Value = "KH041"
$ws.Cells.Item(7, 4).Value = "NV003"
$ws.Cells.Item(7, 6).Value = 0.0

# Row 8
$ws.Cells.Item(8, 1).Value = "HD007"
$ws.Cells.Item(8, 2).Value = "2023-12-09"
$ws.Cells.Item(8, 3).Value = "KH041"
$ws.Cells.Item(8, 4).Value = "NV003"
$ws.Cells.Item(8, 6).Value = 0.0

# Row 9
$ws.Cells.Item(9, 1).Value = "HD008"
$ws.Cells.Item(9, 2).Value = "2023-12-09"
$ws.Cells.Item(9, 3).Value = "KH041"
$ws.Cells.Item(9, 4).Value = "NV003"
$ws.Cells.Item(9, 6).Value = 0.0
